# Fills in manufacturer / reference_number / model_name (columns C, D, E, F)
# for lots 1-6 (rows 2-7), scraped from the listing title that already
# lives in column A. Also (re)creates the other empty "scraped field"
# cells (material, case_number, diameter, movement_number, calibre,
# bracelet_strap, accessoires, signed) for those same rows so the row
# shape matches the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text to put into column C (manufacturer) and F (model_name) - identical
# for every row in this batch - plus D (year) / E (reference_number) where
# they could be parsed out of the title.
$rowData = @{
    2 = @{ C = "NATIONAL WATCH"; F = "NATIONAL WATCH" }
    3 = @{ C = "GALLET"; F = "GALLET" }
    4 = @{ C = "HELVETIA"; F = "HELVETIA" }
    5 = @{
        C = "BELL & ROSS REF. BR 01-97 PVD STEEL LIMITED EDITION`nBell & Ross"
        E = "BR"
        F = "BELL & ROSS REF. BR 01-97 PVD STEEL LIMITED EDITION`nBell & Ross"
    }
    6 = @{
        C = "CHRONOSWISS KLASSIK REF. CH 7443 CHRONOGRAPH STEEL `nChronoswiss"
        D = "7443"
        E = "CH"
        F = "CHRONOSWISS KLASSIK REF. CH 7443 CHRONOGRAPH STEEL `nChronoswiss"
    }
    7 = @{
        C = "VACHERON & CONSTANTIN REF. 33093 YELLOW GOLD`nVacheron & Constantin"
        E = "33093"
        F = "VACHERON & CONSTANTIN REF. 33093 YELLOW GOLD`nVacheron & Constantin"
    }
}

# Every column the scraper writes per row (in sheet order). Any of these
# not given a real value above still gets "touched" so the cell exists
# (empty) in the sheet, matching the other already-populated rows.
$allCols = @("C", "D", "E", "F", "L", "M", "O", "P", "Q", "R", "S", "T")

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    foreach ($col in $allCols) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        # Force text formatting first so purely-numeric-looking strings
        # (e.g. reference numbers like "7443") are stored as text, not
        # auto-converted to numbers - matches the scraped/inline-string
        # source data.
        $cell.NumberFormat = "@"
        if ($values.ContainsKey($col)) {
            $cell.Value = $values[$col]
        } else {
            $cell.Value = ""
        }
        # Put formatting back to the sheet default - only the value
        # (presence of the cell) should change, not its style.
        $cell.Style = "Normal"
    }
}
